$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Replace the hyperlink "https://github.com/gfkl/BiblioJ" with plain
#    (non-hyperlinked) text "https://github.com/gfkl/gghy".
# -----------------------------------------------------------------------

# Remove the hyperlink field itself (this unwraps the <w:hyperlink> but
# keeps the run with its hyperlink-style formatting in place).
if ($d.Hyperlinks.Count -gt 0) {
    $hyperlink = $d.Hyperlinks.Item(1)
    $hyperlink.Delete()
}

# Find the now plain-but-still-styled run that held the old URL text and
# delete its text entirely, so we can retype clean, unstyled text in its
# place.
$oldUrlRange = $d.Content
$found = $oldUrlRange.Find.Execute("https://github.com/gfkl/BiblioJ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $oldUrlRange.Text = ""
}

# Locate the paragraph that used to hold the URL ("URL du projet sur
# Github : ") and append the new, unstyled URL text right after it.
$urlPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*URL du projet sur*") {
        $urlPara = $candidate
    }
}
$insertionPoint = $d.Range($urlPara.Range.End - 1, $urlPara.Range.End - 1)
$insertionPoint.InsertAfter("https://github.com/gfkl/gghy")

# -----------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the trailing empty paragraph to sit
#    right after the "Information sur le projet" heading text (still
#    inside that same paragraph, before its paragraph mark).
# -----------------------------------------------------------------------

$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*Information sur le projet*") {
        $headingPara = $candidate
    }
}
$headingEnd = $headingPara.Range.End - 1

# Re-adding a bookmark with the same name ("_GoBack") automatically
# removes it from wherever it used to be, so the trailing paragraph it
# used to occupy ends up empty again, exactly as in the target document.
#
# NOTE: adding a bookmark directly over a zero-length (collapsed) Range
# placed exactly at the end of a paragraph's text turns out to be
# unreliable in this environment (it silently relocates to the start of
# the document), so a temporary marker character is inserted first to
# anchor a non-collapsed Range; the bookmark is created around that
# character and the character is then removed again, leaving a clean,
# empty bookmark at the desired spot.
$tempAnchor = $d.Range($headingEnd, $headingEnd)
$tempAnchor.InsertAfter("#")

$bookmarkRange = $d.Range($headingEnd, $headingEnd + 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$tempCharRange = $d.Range($headingEnd, $headingEnd + 1)
$tempCharRange.Text = ""
